$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (N_Calib_1=20, N_Calib_2=40)
$ws.Range("C2").Value = -1.492263261150804
$ws.Range("D2").Value = 0.1498323900386123

# Row 3 (N_Calib_1=20, N_Calib_2=60)
$ws.Range("C3").Value = -1.012113203770835
$ws.Range("D3").Value = 0.3224866192301192

# Row 4 (N_Calib_1=20, N_Calib_2=100)
$ws.Range("C4").Value = -1.128493902291183
$ws.Range("D4").Value = 0.2712677721884975

# Row 5 (N_Calib_1=20, N_Calib_2=200)
$ws.Range("C5").Value = 0.4398030073783262
$ws.Range("D5").Value = 0.6643700268620152

# Row 6 (N_Calib_1=40, N_Calib_2=60)
$ws.Range("C6").Value = 0.4664445781751394
$ws.Range("D6").Value = 0.6454831198702604

# Row 7 (N_Calib_1=40, N_Calib_2=100)
$ws.Range("C7").Value = 0.6114228285851228
$ws.Range("D7").Value = 0.5471863621953934

# Row 8 (N_Calib_1=40, N_Calib_2=200)
$ws.Range("C8").Value = 2.009581712069077
$ws.Range("D8").Value = 0.05689811950938362
$ws.Range("G8").Value = "No"

# Row 9 (N_Calib_1=60, N_Calib_2=100)
$ws.Range("C9").Value = -0.008765107634912886
$ws.Range("D9").Value = 0.9930855437655497

# Row 10 (N_Calib_1=60, N_Calib_2=200)
$ws.Range("C10").Value = 1.29444821412549
$ws.Range("D10").Value = 0.2089334315494413

# Row 11 (N_Calib_1=100, N_Calib_2=200)
$ws.Range("C11").Value = 1.357457401657593
$ws.Range("D11").Value = 0.1883988861429227
